$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 698.488914887407
$ws.Range("C2").Value = 123.7749528520503
$ws.Range("D2").Value = 105.2207373988716

$ws.Range("B3").Value = 698.488914887407
$ws.Range("C3").Value = 123.7749528520503
$ws.Range("D3").Value = 105.2207373988716

$ws.Range("B4").Value = 74.45264127021989
$ws.Range("C4").Value = 13.58085457251354
$ws.Range("D4").Value = 31.6244695100989

$ws.Range("B5").Value = 624.036273617187
$ws.Range("C5").Value = 110.1940982795367
$ws.Range("D5").Value = 73.59626788877269

$ws.Range("B6").Value = 772.9415561576268
$ws.Range("C6").Value = 137.3558074245638
$ws.Range("D6").Value = 136.8452069089705
